$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2450.3333
$ws.Range("I98").Value = 2521.1592
$ws.Range("J98").Value = 1671.25
$ws.Range("K98").Value = 2521.1592
$ws.Range("L98").Value = 1671.25
$ws.Range("M98").Value = -1023.1592
$ws.Range("N98").Value = -4667.25
$ws.Range("H113").Value = 3703.6667
$ws.Range("I113").Value = 2052.5
$ws.Range("J113").Value = 7006
$ws.Range("K113").Value = 2052.5
$ws.Range("L113").Value = 7006
$ws.Range("M113").Value = 1201.5
$ws.Range("N113").Value = -13514
$ws.Range("H116").Value = 2993.6667
$ws.Range("I116").Value = 2114.75
$ws.Range("J116").Value = 4751.5
$ws.Range("K116").Value = 2114.75
$ws.Range("L116").Value = 4751.5
$ws.Range("M116").Value = 1327.25
$ws.Range("N116").Value = -11635.5
$ws.Range("H121").Value = 1463.3334
$ws.Range("J121").Value = 1445
$ws.Range("L121").Value = 4335
$ws.Range("N121").Value = -7829
$ws.Range("H122").Value = 2450.3333
$ws.Range("I122").Value = 2521.1592
$ws.Range("J122").Value = 1671.25
$ws.Range("K122").Value = 7563.4776
$ws.Range("L122").Value = 5013.75
$ws.Range("M122").Value = -5113.4776
$ws.Range("N122").Value = -9913.75
$ws.Range("H132").Value = 6809525.5
$ws.Range("I132").Value = 11910613
$ws.Range("K132").Value = 35731839
$ws.Range("M132").Value = -35729309
$ws.Range("H137").Value = 1815.8667
$ws.Range("I137").Value = 1577.5555
$ws.Range("J137").Value = 1918
$ws.Range("K137").Value = 4732.666499999999
$ws.Range("L137").Value = 5754
$ws.Range("M137").Value = -2182.666499999999
$ws.Range("N137").Value = -10854

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 71430376
$ws.Range("I61").Value = 111112750
$ws.Range("J61").Value = 2099.8
$ws.Range("K61").Value = 111112750
$ws.Range("L61").Value = 2099.8
$ws.Range("M61").Value = -111112538
$ws.Range("N61").Value = -2523.8
$ws.Range("H74").Value = 1792.9524
$ws.Range("I74").Value = 1509.0526
$ws.Range("J74").Value = 4490
$ws.Range("K74").Value = 1509.0526
$ws.Range("L74").Value = 4490
$ws.Range("M74").Value = -635.0526
$ws.Range("N74").Value = -6238
$ws.Range("H77").Value = 1792.9524
$ws.Range("I77").Value = 1509.0526
$ws.Range("J77").Value = 4490
$ws.Range("K77").Value = 7545.263
$ws.Range("L77").Value = 22450
$ws.Range("M77").Value = -3177.263
$ws.Range("N77").Value = -31186
$ws.Range("H110").Value = 1328.2858
$ws.Range("I110").Value = 864.17645
$ws.Range("J110").Value = 3300.75
$ws.Range("K110").Value = 864.17645
$ws.Range("L110").Value = 3300.75
$ws.Range("M110").Value = 1180.82355
$ws.Range("N110").Value = -7390.75
$ws.Range("H132").Value = 2783.647
$ws.Range("I132").Value = 2310
$ws.Range("J132").Value = 3652
$ws.Range("K132").Value = 6930
$ws.Range("L132").Value = 10956
$ws.Range("M132").Value = -4400
$ws.Range("N132").Value = -16016
$ws.Range("H136").Value = 71430376
$ws.Range("I136").Value = 111112750
$ws.Range("J136").Value = 2099.8
$ws.Range("K136").Value = 333338250
$ws.Range("L136").Value = 6299.400000000001
$ws.Range("M136").Value = -333335700
$ws.Range("N136").Value = -11399.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4871.893
$ws.Range("I134").Value = 882.4545000000001
$ws.Range("J134").Value = 19499.834
$ws.Range("K134").Value = 2647.3635
$ws.Range("L134").Value = 58499.50199999999
$ws.Range("M134").Value = -112.3635000000004
$ws.Range("N134").Value = -63569.50199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1194.1364
$ws.Range("I31").Value = 1123.2
$ws.Range("J31").Value = 1346.1428
$ws.Range("K31").Value = 1123.2
$ws.Range("L31").Value = 1346.1428
$ws.Range("M31").Value = -828.2
$ws.Range("N31").Value = -1936.1428
$ws.Range("H34").Value = 1194.1364
$ws.Range("I34").Value = 1123.2
$ws.Range("J34").Value = 1346.1428
$ws.Range("K34").Value = 1123.2
$ws.Range("L34").Value = 1346.1428
$ws.Range("M34").Value = -921.2
$ws.Range("N34").Value = -1750.1428
$ws.Range("H58").Value = 1340.2963
$ws.Range("I58").Value = 985
$ws.Range("J58").Value = 2355.4285
$ws.Range("K58").Value = 985
$ws.Range("L58").Value = 2355.4285
$ws.Range("M58").Value = -782
$ws.Range("N58").Value = -2761.4285
$ws.Range("H132").Value = 1712.4595
$ws.Range("I132").Value = 1407.2273
$ws.Range("J132").Value = 2160.1333
$ws.Range("K132").Value = 4221.6819
$ws.Range("L132").Value = 6480.3999
$ws.Range("M132").Value = -1691.6819
$ws.Range("N132").Value = -11540.3999
$ws.Range("H134").Value = 14287462
$ws.Range("I134").Value = 1832.6923
$ws.Range("J134").Value = 55557060
$ws.Range("K134").Value = 5498.0769
$ws.Range("L134").Value = 166671180
$ws.Range("M134").Value = -2963.0769
$ws.Range("N134").Value = -166676250
$ws.Range("H136").Value = 1340.2963
$ws.Range("I136").Value = 985
$ws.Range("J136").Value = 2355.4285
$ws.Range("K136").Value = 2955
$ws.Range("L136").Value = 7066.2855
$ws.Range("M136").Value = -405
$ws.Range("N136").Value = -12166.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 21280204
$ws.Range("I131").Value = 111111420
$ws.Range("J131").Value = 4389.079
$ws.Range("K131").Value = 333334260
$ws.Range("L131").Value = 13167.237
$ws.Range("M131").Value = -333329220
$ws.Range("N131").Value = -23247.237
$ws.Range("H139").Value = 1965.7742
$ws.Range("I139").Value = 2196.7896
$ws.Range("J139").Value = 1600
$ws.Range("K139").Value = 6590.3688
$ws.Range("L139").Value = 4800
$ws.Range("M139").Value = -1450.3688
$ws.Range("N139").Value = -15080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4807.8335
$ws.Range("I132").Value = 5085.25
$ws.Range("J132").Value = 4253
$ws.Range("K132").Value = 15255.75
$ws.Range("L132").Value = 12759
$ws.Range("M132").Value = -12725.75
$ws.Range("N132").Value = -17819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2131.1538
$ws.Range("I7").Value = 1713.25
$ws.Range("J7").Value = 2799.8
$ws.Range("K7").Value = 1713.25
$ws.Range("L7").Value = 2799.8
$ws.Range("M7").Value = -1601.25
$ws.Range("N7").Value = -3023.8
$ws.Range("H40").Value = 2606.4644
$ws.Range("I40").Value = 1847.0435
$ws.Range("K40").Value = 1847.0435
$ws.Range("M40").Value = -1711.0435
$ws.Range("H55").Value = 934.44446
$ws.Range("I55").Value = 801.25
$ws.Range("K55").Value = 801.25
$ws.Range("M55").Value = -628.25
$ws.Range("H126").Value = 2131.1538
$ws.Range("I126").Value = 1713.25
$ws.Range("J126").Value = 2799.8
$ws.Range("K126").Value = 5139.75
$ws.Range("L126").Value = 8399.400000000001
$ws.Range("M126").Value = -2669.75
$ws.Range("N126").Value = -13339.4
$ws.Range("H132").Value = 54870.477
$ws.Range("I132").Value = 14698.375
$ws.Range("K132").Value = 44095.125
$ws.Range("M132").Value = -41565.125
$ws.Range("H136").Value = 1431.5
$ws.Range("I136").Value = 1058.5
$ws.Range("K136").Value = 3175.5
$ws.Range("M136").Value = -625.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1844.862
$ws.Range("I132").Value = 1608.6086
$ws.Range("J132").Value = 2750.5
$ws.Range("K132").Value = 4825.825800000001
$ws.Range("L132").Value = 8251.5
$ws.Range("M132").Value = -2295.825800000001
$ws.Range("N132").Value = -13311.5
$ws.Range("H133").Value = 45715
$ws.Range("J133").Value = 45715
$ws.Range("L133").Value = 45715
$ws.Range("N133").Value = -55835
$ws.Range("H136").Value = 841.1177
$ws.Range("I136").Value = 748.3333
$ws.Range("J136").Value = 1199
$ws.Range("K136").Value = 2244.9999
$ws.Range("L136").Value = 3597
$ws.Range("M136").Value = 305.0001000000002
$ws.Range("N136").Value = -8697
